$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.369.28'
$ws.Range('E2').Value = '  +7.81%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.036.90'
$ws.Range('E3').Value = '  +5.36%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '582.85'
$ws.Range('E5').Value = '  +2.91%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '158.48'
$ws.Range('E6').Value = '  +11.49%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.034.48'
$ws.Range('E8').Value = '  +5.31%  '
$ws.Range('E9').Value = '  +3.85%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.00'
$ws.Range('E10').Value = '  +2.26%  '
$ws.Range('E11').Value = '  +7.66%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.454'
$ws.Range('E12').Value = '  +6.04%  '
$ws.Range('E13').Value = '  +10.24%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.72'
$ws.Range('E14').Value = '  +10.02%  '
$ws.Range('E15').Value = '  +0.81%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.331.60'
$ws.Range('E16').Value = '  +7.83%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.534.87'
$ws.Range('E17').Value = '  +5.21%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.97'
$ws.Range('E18').Value = '  +7.60%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.025.80'
$ws.Range('E19').Value = '  +4.90%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '466.02'
$ws.Range('E20').Value = '  +8.70%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.93'
$ws.Range('E21').Value = '  +7.22%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.687'
$ws.Range('E22').Value = '  +5.84%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.38'
$ws.Range('E23').Value = '  +8.86%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '82.57'
$ws.Range('E24').Value = '  +4.87%  '
$ws.Range('E25').Value = '  +13.06%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.51'
$ws.Range('E26').Value = '  +5.42%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.66'
$ws.Range('E27').Value = '  +6.28%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.07'
$ws.Range('E29').Value = '  +14.94%  '
$ws.Range('E30').Value = '  +16.93%  '
$ws.Range('E31').Value = '  +0.90%  '
$ws.Range('E32').Value = '  +5.15%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.16'
$ws.Range('E33').Value = '  +7.20%  '
$ws.Range('E34').Value = '  +5.90%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +5.24%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.81'
$ws.Range('E37').Value = '  +8.64%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.18'
$ws.Range('E38').Value = '  +14.95%  '
$ws.Range('E39').Value = '  +10.45%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '49.66'
$ws.Range('E40').Value = '  +1.84%  '
$ws.Range('E41').Value = '  +8.29%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '44.02'
$ws.Range('E42').Value = '  +12.47%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.303'
$ws.Range('E43').Value = '  +14.23%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.46'
$ws.Range('E44').Value = '  +4.07%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '391.35'
$ws.Range('E45').Value = '  +14.05%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.813.45'
$ws.Range('E46').Value = '  +5.10%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0355'
$ws.Range('E47').Value = '  +6.69%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '134.24'
$ws.Range('E48').Value = '  +1.67%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '23.64'
$ws.Range('E50').Value = '  +10.70%  '
$ws.Range('E51').Value = '  +4.83%  '
